$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4286.5713
$ws.Range("I19").Value = 10487
$ws.Range("J19").Value = 841.8889
$ws.Range("K19").Value = 10487
$ws.Range("L19").Value = 841.8889
$ws.Range("M19").Value = -10312
$ws.Range("N19").Value = -1191.8889
$ws.Range("H80").Value = 1575.4546
$ws.Range("I80").Value = 573.3333
$ws.Range("J80").Value = 1951.25
$ws.Range("K80").Value = 1719.9999
$ws.Range("L80").Value = 5853.75
$ws.Range("M80").Value = -721.9999
$ws.Range("N80").Value = -7849.75
$ws.Range("H83").Value = 1575.4546
$ws.Range("I83").Value = 573.3333
$ws.Range("J83").Value = 1951.25
$ws.Range("K83").Value = 5159.9997
$ws.Range("L83").Value = 17561.25
$ws.Range("M83").Value = -167.9997000000003
$ws.Range("N83").Value = -27545.25
$ws.Range("H129").Value = 1354.7727
$ws.Range("I129").Value = 474.2857
$ws.Range("J129").Value = 1765.6666
$ws.Range("K129").Value = 1422.8571
$ws.Range("L129").Value = 5296.9998
$ws.Range("M129").Value = 3577.1429
$ws.Range("N129").Value = -15296.9998
$ws.Range("H138").Value = 14288104
$ws.Range("I138").Value = 1123.5143
$ws.Range("J138").Value = 28575084
$ws.Range("K138").Value = 3370.5429
$ws.Range("L138").Value = 85725252
$ws.Range("M138").Value = 1769.4571
$ws.Range("N138").Value = -85735532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 901.6957
$ws.Range("I2").Value = 850.7368
$ws.Range("J2").Value = 1143.75
$ws.Range("K2").Value = 850.7368
$ws.Range("L2").Value = 1143.75
$ws.Range("M2").Value = -737.7368
$ws.Range("N2").Value = -1369.75
$ws.Range("H32").Value = 5757.4614
$ws.Range("I32").Value = 3258.6924
$ws.Range("J32").Value = 13253.77
$ws.Range("K32").Value = 3258.6924
$ws.Range("L32").Value = 13253.77
$ws.Range("M32").Value = -2971.6924
$ws.Range("N32").Value = -13827.77
$ws.Range("H45").Value = 1946.4117
$ws.Range("I45").Value = 1016
$ws.Range("J45").Value = 2597.7
$ws.Range("K45").Value = 1016
$ws.Range("L45").Value = 2597.7
$ws.Range("M45").Value = -639
$ws.Range("N45").Value = -3351.7
$ws.Range("H113").Value = 38966.168
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 38966.168
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 38966.168
$ws.Range("N113").Value = -47644.168
$ws.Range("H116").Value = 901.6957
$ws.Range("I116").Value = 850.7368
$ws.Range("J116").Value = 1143.75
$ws.Range("K116").Value = 850.7368
$ws.Range("L116").Value = 1143.75
$ws.Range("M116").Value = 1443.2632
$ws.Range("N116").Value = -5731.75
$ws.Range("H139").Value = 91073
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 91073
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 91073
$ws.Range("N139").Value = -101353

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 901.6957
$ws.Range("I3").Value = 850.7368
$ws.Range("J3").Value = 1143.75
$ws.Range("K3").Value = 850.7368
$ws.Range("L3").Value = 1143.75
$ws.Range("M3").Value = -736.7368
$ws.Range("N3").Value = -1371.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17246286
$ws.Range("I31").Value = 38462400
$ws.Range("J31").Value = 8192.718999999999
$ws.Range("K31").Value = 38462400
$ws.Range("L31").Value = 8192.718999999999
$ws.Range("M31").Value = -38462105
$ws.Range("N31").Value = -8782.718999999999
$ws.Range("H34").Value = 17246286
$ws.Range("I34").Value = 38462400
$ws.Range("J34").Value = 8192.718999999999
$ws.Range("K34").Value = 38462400
$ws.Range("L34").Value = 8192.718999999999
$ws.Range("M34").Value = -38462198
$ws.Range("N34").Value = -8596.718999999999
$ws.Range("H58").Value = 2656
$ws.Range("I58").Value = 5506
$ws.Range("J58").Value = 1516
$ws.Range("K58").Value = 5506
$ws.Range("L58").Value = 1516
$ws.Range("M58").Value = -5303
$ws.Range("N58").Value = -1922
$ws.Range("H98").Value = 78624
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 78624
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 78624
$ws.Range("N98").Value = -83116
$ws.Range("H124").Value = 18724.75
$ws.Range("I124").Value = 20000
$ws.Range("J124").Value = 18542.572
$ws.Range("K124").Value = 20000
$ws.Range("L124").Value = 18542.572
$ws.Range("M124").Value = -17545
$ws.Range("N124").Value = -23452.572
$ws.Range("H134").Value = 3121.4565
$ws.Range("I134").Value = 4364.6924
$ws.Range("J134").Value = 1505.25
$ws.Range("K134").Value = 13094.0772
$ws.Range("L134").Value = 4515.75
$ws.Range("M134").Value = -10559.0772
$ws.Range("N134").Value = -9585.75
$ws.Range("H136").Value = 2656
$ws.Range("I136").Value = 5506
$ws.Range("J136").Value = 1516
$ws.Range("K136").Value = 16518
$ws.Range("L136").Value = 4548
$ws.Range("M136").Value = -13968
$ws.Range("N136").Value = -9648

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1590.5
$ws.Range("I34").Value = 116.666664
$ws.Range("J34").Value = 1930.6154
$ws.Range("K34").Value = 349.999992
$ws.Range("L34").Value = 5791.8462
$ws.Range("M34").Value = -265.999992
$ws.Range("N34").Value = -5959.8462
$ws.Range("H68").Value = 1239388.5
$ws.Range("I68").Value = 2431796.5
$ws.Range("J68").Value = 1118.5
$ws.Range("K68").Value = 7295389.5
$ws.Range("L68").Value = 3355.5
$ws.Range("M68").Value = -7294578.5
$ws.Range("N68").Value = -4977.5
$ws.Range("H71").Value = 1239388.5
$ws.Range("I71").Value = 2431796.5
$ws.Range("J71").Value = 1118.5
$ws.Range("K71").Value = 21886168.5
$ws.Range("L71").Value = 10066.5
$ws.Range("M71").Value = -21882112.5
$ws.Range("N71").Value = -18178.5
$ws.Range("H107").Value = 798.46344
$ws.Range("I107").Value = 832.04346
$ws.Range("J107").Value = 755.55554
$ws.Range("K107").Value = 2496.13038
$ws.Range("L107").Value = 2266.66662
$ws.Range("M107").Value = -576.1303800000001
$ws.Range("N107").Value = -6106.66662
$ws.Range("H131").Value = 1097.15
$ws.Range("I131").Value = 461.125
$ws.Range("J131").Value = 1152.4565
$ws.Range("K131").Value = 1383.375
$ws.Range("L131").Value = 3457.3695
$ws.Range("M131").Value = 3656.625
$ws.Range("N131").Value = -13537.3695

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11865.615
$ws.Range("I70").Value = 15088.667
$ws.Range("J70").Value = 4613.75
$ws.Range("K70").Value = 15088.667
$ws.Range("L70").Value = 4613.75
$ws.Range("M70").Value = -14818.667
$ws.Range("N70").Value = -5153.75
$ws.Range("H73").Value = 11865.615
$ws.Range("I73").Value = 15088.667
$ws.Range("J73").Value = 4613.75
$ws.Range("K73").Value = 15088.667
$ws.Range("L73").Value = 4613.75
$ws.Range("M73").Value = -14152.667
$ws.Range("N73").Value = -6485.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1939.9333
$ws.Range("I100").Value = 1930.6923
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1930.6923
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1389.6923
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 10843.667
$ws.Range("I88").Value = 1171
$ws.Range("J88").Value = 30189
$ws.Range("K88").Value = 1171
$ws.Range("L88").Value = 30189
$ws.Range("M88").Value = -765
$ws.Range("N88").Value = -31001
$ws.Range("H91").Value = 10843.667
$ws.Range("I91").Value = 1171
$ws.Range("J91").Value = 30189
$ws.Range("K91").Value = 1171
$ws.Range("L91").Value = 30189
$ws.Range("M91").Value = 233
$ws.Range("N91").Value = -32997
$ws.Range("H122").Value = 68422.266
$ws.Range("I122").Value = 126542.375
$ws.Range("J122").Value = 1999.2858
$ws.Range("K122").Value = 379627.125
$ws.Range("L122").Value = 5997.857400000001
$ws.Range("M122").Value = -377177.125
$ws.Range("N122").Value = -10897.8574
$ws.Range("H136").Value = 4079.4146
$ws.Range("I136").Value = 647.8
$ws.Range("J136").Value = 7347.619
$ws.Range("K136").Value = 1943.4
$ws.Range("L136").Value = 22042.857
$ws.Range("M136").Value = 606.6000000000001
$ws.Range("N136").Value = -27142.857
